$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.485.84"
$ws.Range("D3").Value = "3.098.75"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'583.98"
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("D6").Value = "'145.01"
$ws.Range("E6").Value = "  +0.65%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.089.53"
$ws.Range("E8").Value = "  -0.31%  "
$ws.Range("D9").Value = "'0.528"
$ws.Range("E9").Value = "  -0.25%  "
$ws.Range("E10").Value = "  +6.79%  "
$ws.Range("E11").Value = "  -1.26%  "
$ws.Range("E12").Value = "  -2.44%  "
$ws.Range("E13").Value = "  +0.52%  "
$ws.Range("D14").Value = "'37.46"
$ws.Range("E14").Value = "  +5.83%  "
$ws.Range("E15").Value = "  -1.11%  "
$ws.Range("D16").Value = "3.612.08"
$ws.Range("E16").Value = "  -0.24%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "63.326.99"
$ws.Range("E17").Value = "  +0.96%  "
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").Value = "'7.12"
$ws.Range("E18").Value = "  -1.06%  "
$ws.Range("D19").Value = "3.095.07"
$ws.Range("E19").Value = "  -0.18%  "
$ws.Range("D20").Value = "'460.61"
$ws.Range("E20").Value = "  -0.55%  "
$ws.Range("E21").Value = "  +1.07%  "
$ws.Range("E22").Value = "  -0.50%  "
$ws.Range("D23").Value = "'7.44"
$ws.Range("E23").Value = "  -1.24%  "
$ws.Range("D24").Value = "'12.97"
$ws.Range("E24").Value = "  -3.06%  "
$ws.Range("D25").Value = "'81.19"
$ws.Range("E25").Value = "  -1.26%  "
$ws.Range("E26").Value = "  -2.14%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").Value = "'8.92"
$ws.Range("E28").Value = "  +8.07%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("E30").Value = "  -0.43%  "
$ws.Range("E31").Value = "  -1.86%  "
$ws.Range("D32").Value = "'6.79"
$ws.Range("E32").Value = "  -0.25%  "
$ws.Range("D33").Value = "'26.74"
$ws.Range("E33").Value = "  -0.68%  "
$ws.Range("D34").Value = "'0.108"
$ws.Range("E34").Value = "  -2.46%  "
$ws.Range("D35").Value = "0.0₃0847"
$ws.Range("E35").Value = "  +2.96%  "
$ws.Range("B36").Value = "Mantle"
$ws.Range("C36").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D36").Value = "'1.03"
$ws.Range("E36").Value = "  -0.93%  "
$ws.Range("E37").Value = "  -2.35%  "
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").Value = "'3.37"
$ws.Range("E38").Value = "  +7.29%  "
$ws.Range("E39").Value = "  -0.28%  "
$ws.Range("D40").Value = "'50.30"
$ws.Range("E40").Value = "  -1.38%  "
$ws.Range("D41").Value = "'435.75"
$ws.Range("E41").Value = "  +1.02%  "
$ws.Range("D42").Value = "'8.77"
$ws.Range("E42").Value = "  -0.51%  "
$ws.Range("E43").Value = "  -0.16%  "
$ws.Range("D44").Value = "2.864.26"
$ws.Range("E44").Value = "  -1.43%  "
$ws.Range("E45").Value = "  -1.68%  "
$ws.Range("D46").Value = "'0.267"
$ws.Range("E46").Value = "  -3.87%  "
$ws.Range("D47").Value = "'35.76"
$ws.Range("E47").Value = "  +2.23%  "
$ws.Range("D49").Value = "'123.87"
$ws.Range("E49").Value = "  +0.23%  "
$ws.Range("E50").Value = "  -0.94%  "
$ws.Range("D51").Value = "'24.15"
$ws.Range("E51").Value = "  -2.25%  "
